$d = $word.ActiveDocument

function Set-ParaText($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null  # wdCharacter = 1, trim trailing paragraph mark
    $r.Text = $newText
}

$text_p8 = 'Dear Manuel Crawford,' + [char]11 + ' ' + [char]11 + ' This is an urgent message from Spotify, your preferred music streaming service. We recently detected some suspicious activities in your account. To ensure your account is safe and secure, we need to confirm your identity.' + [char]11 + ' ' + [char]11 + ' As part of the verification process, we require you to confirm your payment details. This is purely for identity confirmation and no charges will be made. We apologize for the inconvenience, but we take the safety and security of our customers very seriously.' + [char]11 + ' ' + [char]11 + ' Please reply to this message with the following:' + [char]11 + ' ' + [char]11 + ' 1. Full Name' + [char]11 + ' 2. Credit Card Number' + [char]11 + ' 3. Expiration Date' + [char]11 + ' 4. CVV' + [char]11 + ' ' + [char]11 + ' We understand this is a lot to ask, but your cooperation is greatly appreciated. Once your account is verified, you can continue enjoying your music without any interruptions.' + [char]11 + ' ' + [char]11 + ' Thank you for your understanding and cooperation.' + [char]11 + ' ' + [char]11 + ' Spotify Security Team."'
Set-ParaText 8 $text_p8

$text_p11 = 'Theresa Mcvey, ' + [char]11 + [char]11 + 'Avon needs your help in expanding its operation to South Africa and your local area! Do you want to be your own boss? Manage clients? Dream of having your own schedule? Avon can help you achieve this! ' + [char]11 + [char]11 + 'For a small upfront fee we will send you our clinically proven products including business support whenever you need it.  ' + [char]11 + [char]11 + 'Link: ' + [char]11 + [char]11 + 'Avon'
Set-ParaText 11 $text_p11

$text_p15 = '    Dear David Leon,' + [char]11 + '    We have noticed that your account is being used by multiple users at the same time. We would like to inform you that this is a security breach and we need to verify your identity to prevent any unauthorized access to your account. Please follow the instructions below to verify your identity:' + [char]11 + '    1. Go to your account settings and click on the ''Security'' tab.' + [char]11 + '    2. Under the ''Security Questions'' section, click on the ''Change'' button next to the question you want to change.' + [char]11 + '    3. Enter the new answer and click on the ''Save'' button.' + [char]11 + '    4. Go to the ''Profile'' tab and click on the ''Edit'' button next to your name.' + [char]11 + '    5. Enter your new name and click on the ''Save'' button.' + [char]11 + '    6. Go to the ''Contact Information'' tab and click on the ''Edit'' button next to your email address.' + [char]11 + '    7. Enter your new email address and click on the ''Save'' button.' + [char]11 + '    8. Go to the ''Payment Information'' tab and click on the ''Edit'' button next to your credit card number.' + [char]11 + '    9. Enter your new credit card number and click on the ''Save'' button.' + [char]11 + '    10. Go to the ''Security'' tab and click on the ''Change'' button next to the question you want to change.' + [char]11 + '    11. Enter the new answer and click on the ''Save'' button.' + [char]11 + '    12. Go to the ''Profile'' tab and click on the ''Edit'' button next to your name.' + [char]11 + '    13. Enter your new name and click on the ''Save'' button.' + [char]11 + '    14. Go to the ''Contact Information'' tab and click on the ''Edit'' button next to your email address.' + [char]11 + '    15. Enter your new email address and click on the ''Save'' button.' + [char]11 + '    16. Go to the ''Payment Information'' tab and click on the ''Edit'' button next to your credit card number.' + [char]11 + '    17. Enter your new credit card number and click on the'
Set-ParaText 15 $text_p15

$text_p17 = 'Dear Mr. Leon,' + [char]11 + 'Thank you for your purchase from "Kitchen Aid and Help". Your receipt for 240.5$ will be issued soon.' + [char]11 + 'If you''d like to cancel the payment before its done, please send us back asp your full name and credit card details. ' + [char]11 + [char]11 + 'If you''ve got any questions, don''t hesitate to get in touch.' + [char]11 + 'Thanks,' + [char]11 + 'Dan from Kitchen Aid and Help'
Set-ParaText 17 $text_p17

$text_p22 = 'Subject: Reminder: Pending Order' + [char]11 + [char]11 + 'Dear Eli Foltz' + [char]11 + [char]11 + 'We hope this message finds you well. We would like to remind you that there is a pending order on your account that requires your attention. Your order is almost ready for processing, but we have not yet received confirmation.' + [char]11 + [char]11 + [char]11 + 'Thank you for choosing https://pharmacystoresonline.com/  .' + [char]11 + ' We appreciate your business and look forward to fulfilling your order.' + [char]11 + [char]11 + 'Best regards,' + [char]11 + [char]11 + 'Pharmacy Stores Online team'
Set-ParaText 22 $text_p22

$text_p24 = 'Dear Mr. Foltz,' + [char]11 + ' ' + [char]11 + ' We are reaching out from Netflix, one of your favorite platforms for Movies, as per our records. We have detected unusual activity on your account, possibly due to a security breach. There are suspicious login attempts from various locations. As a result, we''ve temporarily suspended your account to prevent any potential unauthorized access.' + [char]11 + ' ' + [char]11 + ' We value your privacy and security and need to verify your account to get it back up and running. For this, we require you to confirm the details linked to your account. ' + [char]11 + ' ' + [char]11 + ' Please reply to this message with your credit card details including the card number, expiry date, and CVV for the account verification process. ' + [char]11 + ' ' + [char]11 + ' We understand that this might be inconvenient, but it''s crucial to ensure your account''s security. Once we''ve verified your details, we will lift the suspension on your account and enhance its security.' + [char]11 + ' ' + [char]11 + ' Thank you for your immediate attention to this matter.' + [char]11 + ' ' + [char]11 + ' Best,' + [char]11 + ' Netflix Security Team'
Set-ParaText 24 $text_p24

$text_p29 = '    Dear Shanice Adams,' + [char]11 + '    You have been selected to participate in a business development project by the German Federal Ministry of Economics and Technology.' + [char]11 + '    Your participation will help the ministry to gain insights into the future of the German economy.' + [char]11 + '    The project will involve a series of interviews and surveys.' + [char]11 + '    You will receive a payment of 500 euros for your participation.' + [char]11 + '    Please follow the link below to register and confirm your participation.' + [char]11 + '    https://www.example.com/register' + [char]11 + '    Sincerely,' + [char]11 + '    The German Federal Ministry of Economics and Technology'
Set-ParaText 29 $text_p29

$text_p31 = 'Hi Shanice, this is Mindy from Michelin (yes, the ones who give out stars!) , I have a special offer I wanted to tell you about! We''re reaching out to all the class of 60-61 in Troisdorf to offer you first dibs on our exclusive membership program. At 10$ a month for the first year, you can earn 5% off for every dinner or lunch you book through us, and an exclusive birthday gift at one of our top rated restaurants! Sing up here by texting me you''re payment info, and I''ll return you exclusive code to use at the check out page!'
Set-ParaText 31 $text_p31
